$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 90808
$ws.Range("I62").Value = 161400.58
$ws.Range("J62").Value = 8450
$ws.Range("K62").Value = 161400.58
$ws.Range("L62").Value = 8450
$ws.Range("M62").Value = -160776.58
$ws.Range("N62").Value = -9698
$ws.Range("H65").Value = 90808
$ws.Range("I65").Value = 161400.58
$ws.Range("J65").Value = 8450
$ws.Range("K65").Value = 807002.8999999999
$ws.Range("L65").Value = 42250
$ws.Range("M65").Value = -803882.8999999999
$ws.Range("N65").Value = -48490
$ws.Range("H76").Value = 6252937.5
$ws.Range("I76").Value = 6252937.5
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 6252937.5
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -6252622.5
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 6252937.5
$ws.Range("I79").Value = 6252937.5
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 6252937.5
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -6251845.5
$ws.Range("N79").ClearContents()
$ws.Range("H86").Value = 41670252
$ws.Range("I86").Value = 2809.182
$ws.Range("K86").Value = 2809.182
$ws.Range("M86").Value = -1686.182
$ws.Range("H89").Value = 41670252
$ws.Range("I89").Value = 2809.182
$ws.Range("K89").Value = 14045.91
$ws.Range("M89").Value = -8429.91
$ws.Range("H96").Value = 897
$ws.Range("I96").Value = 1062.7778
$ws.Range("J96").Value = 399.66666
$ws.Range("K96").Value = 3188.3334
$ws.Range("L96").Value = 1198.99998
$ws.Range("M96").Value = -1815.3334
$ws.Range("N96").Value = -3944.99998
$ws.Range("H98").Value = 1481.8182
$ws.Range("I98").Value = 1611.1111
$ws.Range("J98").Value = 900
$ws.Range("K98").Value = 1611.1111
$ws.Range("L98").Value = 900
$ws.Range("M98").Value = -113.1111000000001
$ws.Range("N98").Value = -3896
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H113").Value = 2996.8845
$ws.Range("I113").Value = 2951.389
$ws.Range("J113").Value = 3099.25
$ws.Range("K113").Value = 2951.389
$ws.Range("L113").Value = 3099.25
$ws.Range("M113").Value = 302.6109999999999
$ws.Range("N113").Value = -9607.25
$ws.Range("H116").Value = 4441.875
$ws.Range("I116").Value = 4810
$ws.Range("J116").Value = 4155.5557
$ws.Range("K116").Value = 4810
$ws.Range("L116").Value = 4155.5557
$ws.Range("M116").Value = -1368
$ws.Range("N116").Value = -11039.5557
$ws.Range("H122").Value = 1481.8182
$ws.Range("I122").Value = 1611.1111
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 4833.3333
$ws.Range("L122").Value = 2700
$ws.Range("M122").Value = -2383.3333
$ws.Range("N122").Value = -7600
$ws.Range("H132").Value = 4940.5
$ws.Range("I132").Value = 5156.1113
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 15468.3339
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -12938.3339
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2278.037
$ws.Range("I2").Value = 1878.7222
$ws.Range("J2").Value = 3076.6667
$ws.Range("K2").Value = 1878.7222
$ws.Range("L2").Value = 3076.6667
$ws.Range("M2").Value = -1765.7222
$ws.Range("N2").Value = -3302.6667
$ws.Range("H61").Value = 2500
$ws.Range("I61").Value = 2000
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 2000
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -1788
$ws.Range("N61").Value = -3424
$ws.Range("H74").Value = 1125.6818
$ws.Range("I74").Value = 1077.55
$ws.Range("K74").Value = 1077.55
$ws.Range("M74").Value = -203.55
$ws.Range("H77").Value = 1125.6818
$ws.Range("I77").Value = 1077.55
$ws.Range("K77").Value = 5387.75
$ws.Range("M77").Value = -1019.75
$ws.Range("H116").Value = 2278.037
$ws.Range("I116").Value = 1878.7222
$ws.Range("J116").Value = 3076.6667
$ws.Range("K116").Value = 1878.7222
$ws.Range("L116").Value = 3076.6667
$ws.Range("M116").Value = 415.2778000000001
$ws.Range("N116").Value = -7664.6667
$ws.Range("H136").Value = 2500
$ws.Range("I136").Value = 2000
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -3450
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2278.037
$ws.Range("I3").Value = 1878.7222
$ws.Range("J3").Value = 3076.6667
$ws.Range("K3").Value = 1878.7222
$ws.Range("L3").Value = 3076.6667
$ws.Range("M3").Value = -1764.7222
$ws.Range("N3").Value = -3304.6667
$ws.Range("H140").Value = 76212.86
$ws.Range("J140").Value = 76212.86
$ws.Range("L140").Value = 76212.86
$ws.Range("N140").Value = -86572.86

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I16").Value = 1600
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1600
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1313
$ws.Range("N16").ClearContents()
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("I113").Value = 1600
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1600
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 570
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5895135
$ws.Range("J131").Value = 958.9103
$ws.Range("L131").Value = 2876.7309
$ws.Range("N131").Value = -12956.7309

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H70").Value = 4262.5
$ws.Range("I70").Value = 4257.143
$ws.Range("J70").Value = 4300
$ws.Range("K70").Value = 4257.143
$ws.Range("L70").Value = 4300
$ws.Range("M70").Value = -3987.143
$ws.Range("N70").Value = -4840
$ws.Range("H73").Value = 4262.5
$ws.Range("I73").Value = 4257.143
$ws.Range("J73").Value = 4300
$ws.Range("K73").Value = 4257.143
$ws.Range("L73").Value = 4300
$ws.Range("M73").Value = -3321.143
$ws.Range("N73").Value = -6172
$ws.Range("H140").Value = 65040
$ws.Range("J140").Value = 65040
$ws.Range("L140").Value = 65040
$ws.Range("N140").Value = -75400

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 10203437
$ws.Range("I100").Value = 14029089
$ws.Range("J100").Value = 1700
$ws.Range("K100").Value = 14029089
$ws.Range("L100").Value = 1700
$ws.Range("M100").Value = -14028548
$ws.Range("N100").Value = -2782
